$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: QuestItem / QuestItemCount (F3/G3) ---
# F3 10202002 -> 10114101, and its style goes from centered (s="1") to general (s default)
$ws.Range("F3").Value2 = 10114101
$ws.Range("F3").HorizontalAlignment = 1   # xlGeneral
$ws.Range("G3").Value2 = 1

# --- Row 5: QuestMonster2 / QuestMonster2Count (L5/M5) cleared, style kept ---
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()

# --- Row 10: QuestMonster (J10) 9999 -> 1 ---
$ws.Range("J10").Value2 = 1

# --- Rows 15-18, 20: add QuestMonster / QuestMonsterCount (J/K) = 1 / 5 ---
$ws.Range("J15").Value2 = 1
$ws.Range("K15").Value2 = 5

$ws.Range("J16").Value2 = 1
$ws.Range("K16").Value2 = 5

$ws.Range("J17").Value2 = 1
$ws.Range("K17").Value2 = 5

$ws.Range("J18").Value2 = 1
$ws.Range("K18").Value2 = 5

$ws.Range("J20").Value2 = 1
$ws.Range("K20").Value2 = 5

# --- Selection moved from F4 to M5 ---
$ws.Range("M5").Select()
